$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3710, 4022, 4481, 4540, 4722, 4722, 4805, 4805, 4805, 4916, 4972, 4972, 4989, 4989)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
